$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 312.75
$ws.Cells.Item(2, 9).Value = 274.57144
$ws.Cells.Item(2, 10).Value = 580
$ws.Cells.Item(2, 11).Value = 274.57144
$ws.Cells.Item(2, 12).Value = 580
$ws.Cells.Item(2, 13).Value = -161.57144
$ws.Cells.Item(2, 14).Value = -806
$ws.Cells.Item(15, 8).Value = 1101.2976
$ws.Cells.Item(15, 9).Value = 1101.2976
$ws.Cells.Item(15, 11).Value = 3303.892800000001
$ws.Cells.Item(15, 13).Value = -3134.892800000001
$ws.Cells.Item(38, 8).Value = 318.375
$ws.Cells.Item(38, 9).Value = 336.7143
$ws.Cells.Item(38, 10).Value = 190
$ws.Cells.Item(38, 11).Value = 1010.1429
$ws.Cells.Item(38, 12).Value = 570
$ws.Cells.Item(38, 13).Value = -638.1428999999999
$ws.Cells.Item(38, 14).Value = -1314
$ws.Cells.Item(51, 8).Value = 18280.578
$ws.Cells.Item(51, 10).Value = 6364.154
$ws.Cells.Item(51, 12).Value = 6364.154
$ws.Cells.Item(51, 14).Value = -7332.154
$ws.Cells.Item(58, 8).Value = 8808.706
$ws.Cells.Item(58, 10).Value = 14518.7
$ws.Cells.Item(58, 12).Value = 43556.10000000001
$ws.Cells.Item(58, 14).Value = -43856.10000000001
$ws.Cells.Item(64, 8).Value = 5561.3335
$ws.Cells.Item(64, 9).Value = 4343.8335
$ws.Cells.Item(64, 10).Value = 7996.3335
$ws.Cells.Item(64, 11).Value = 4343.8335
$ws.Cells.Item(64, 12).Value = 7996.3335
$ws.Cells.Item(64, 13).Value = -4095.8335
$ws.Cells.Item(64, 14).Value = -8492.333500000001
$ws.Cells.Item(67, 8).Value = 5561.3335
$ws.Cells.Item(67, 9).Value = 4343.8335
$ws.Cells.Item(67, 10).Value = 7996.3335
$ws.Cells.Item(67, 11).Value = 4343.8335
$ws.Cells.Item(67, 12).Value = 7996.3335
$ws.Cells.Item(67, 13).Value = -3485.8335
$ws.Cells.Item(67, 14).Value = -9712.333500000001
$ws.Cells.Item(82, 8).Value = 6951.143
$ws.Cells.Item(82, 9).Value = 1332.2
$ws.Cells.Item(82, 11).Value = 3996.6
$ws.Cells.Item(82, 13).Value = -3590.6
$ws.Cells.Item(85, 8).Value = 6951.143
$ws.Cells.Item(85, 9).Value = 1332.2
$ws.Cells.Item(85, 11).Value = 3996.6
$ws.Cells.Item(85, 13).Value = -2592.6
$ws.Cells.Item(135, 8).Value = 1756.2
$ws.Cells.Item(135, 9).Value = 1436.4
$ws.Cells.Item(135, 11).Value = 12927.6
$ws.Cells.Item(135, 13).Value = -10392.6
$ws.Cells.Item(138, 8).Value = 6318.9775
$ws.Cells.Item(138, 10).Value = 5343.968
$ws.Cells.Item(138, 12).Value = 16031.904
$ws.Cells.Item(138, 14).Value = -26311.904
$ws.Cells.Item(141, 8).Value = 4956.5
$ws.Cells.Item(141, 9).Value = 2679.9092
$ws.Cells.Item(141, 10).Value = 29999
$ws.Cells.Item(141, 11).Value = 8039.7276
$ws.Cells.Item(141, 12).Value = 89997
$ws.Cells.Item(141, 13).Value = -2859.7276
$ws.Cells.Item(141, 14).Value = -100357

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1930.9642
$ws.Cells.Item(74, 9).Value = 1552.7826
$ws.Cells.Item(74, 11).Value = 1552.7826
$ws.Cells.Item(74, 13).Value = -678.7826
$ws.Cells.Item(77, 8).Value = 1930.9642
$ws.Cells.Item(77, 9).Value = 1552.7826
$ws.Cells.Item(77, 11).Value = 7763.913
$ws.Cells.Item(77, 13).Value = -3395.913

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1479.1666
$ws.Cells.Item(64, 9).Value = 1272
$ws.Cells.Item(64, 10).Value = 1893.5
$ws.Cells.Item(64, 11).Value = 1272
$ws.Cells.Item(64, 12).Value = 1893.5
$ws.Cells.Item(64, 13).Value = -1047
$ws.Cells.Item(64, 14).Value = -2343.5
$ws.Cells.Item(67, 8).Value = 1479.1666
$ws.Cells.Item(67, 9).Value = 1272
$ws.Cells.Item(67, 10).Value = 1893.5
$ws.Cells.Item(67, 11).Value = 1272
$ws.Cells.Item(67, 12).Value = 1893.5
$ws.Cells.Item(67, 13).Value = -492
$ws.Cells.Item(67, 14).Value = -3453.5
$ws.Cells.Item(94, 8).Value = 14549.5
$ws.Cells.Item(94, 9).Value = 279.5
$ws.Cells.Item(94, 11).Value = 279.5
$ws.Cells.Item(94, 13).Value = 171.5
$ws.Cells.Item(105, 8).Value = 3036.1875
$ws.Cells.Item(105, 9).Value = 1980.091
$ws.Cells.Item(105, 11).Value = 1980.091
$ws.Cells.Item(105, 13).Value = -233.0909999999999
$ws.Cells.Item(134, 8).Value = 4045.2954
$ws.Cells.Item(134, 9).Value = 3493.295
$ws.Cells.Item(134, 10).Value = 8350.9
$ws.Cells.Item(134, 11).Value = 10479.885
$ws.Cells.Item(134, 12).Value = 25052.7
$ws.Cells.Item(134, 13).Value = -7944.885
$ws.Cells.Item(134, 14).Value = -30122.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 29734.395
$ws.Cells.Item(31, 9).Value = 3138.147
$ws.Cells.Item(31, 10).Value = 255802.5
$ws.Cells.Item(31, 11).Value = 3138.147
$ws.Cells.Item(31, 12).Value = 255802.5
$ws.Cells.Item(31, 13).Value = -2843.147
$ws.Cells.Item(31, 14).Value = -256392.5
$ws.Cells.Item(34, 8).Value = 29734.395
$ws.Cells.Item(34, 9).Value = 3138.147
$ws.Cells.Item(34, 10).Value = 255802.5
$ws.Cells.Item(34, 11).Value = 3138.147
$ws.Cells.Item(34, 12).Value = 255802.5
$ws.Cells.Item(34, 13).Value = -2936.147
$ws.Cells.Item(34, 14).Value = -256206.5
$ws.Cells.Item(62, 8).Value = 4201.533
$ws.Cells.Item(62, 9).Value = 3856.25
$ws.Cells.Item(62, 11).Value = 3856.25
$ws.Cells.Item(62, 13).Value = -3232.25
$ws.Cells.Item(65, 8).Value = 4201.533
$ws.Cells.Item(65, 9).Value = 3856.25
$ws.Cells.Item(65, 11).Value = 19281.25
$ws.Cells.Item(65, 13).Value = -16161.25
$ws.Cells.Item(88, 8).Value = 60000
$ws.Cells.Item(88, 10).Value = 60000
$ws.Cells.Item(88, 12).Value = 60000
$ws.Cells.Item(88, 14).Value = -60812
$ws.Cells.Item(91, 8).Value = 60000
$ws.Cells.Item(91, 10).Value = 60000
$ws.Cells.Item(91, 12).Value = 60000
$ws.Cells.Item(91, 14).Value = -62808
$ws.Cells.Item(99, 8).Value = 2449.8
$ws.Cells.Item(99, 9).Value = 2449.8
$ws.Cells.Item(99, 11).Value = 2449.8
$ws.Cells.Item(99, 13).Value = -951.8000000000002
$ws.Cells.Item(106, 8).Value = 10200
$ws.Cells.Item(106, 10).Value = 19000
$ws.Cells.Item(106, 12).Value = 19000
$ws.Cells.Item(106, 14).Value = -21524
$ws.Cells.Item(126, 8).Value = 2449.8
$ws.Cells.Item(126, 9).Value = 2449.8
$ws.Cells.Item(126, 11).Value = 7349.400000000001
$ws.Cells.Item(126, 13).Value = -4879.400000000001
$ws.Cells.Item(141, 8).Value = 582999.2
$ws.Cells.Item(141, 10).Value = 685599
$ws.Cells.Item(141, 12).Value = 685599
$ws.Cells.Item(141, 14).Value = -695959

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 392.66666
$ws.Cells.Item(86, 9).Value = 409.16666
$ws.Cells.Item(86, 11).Value = 1227.49998
$ws.Cells.Item(86, 13).Value = -41.49998000000005
$ws.Cells.Item(89, 8).Value = 392.66666
$ws.Cells.Item(89, 9).Value = 409.16666
$ws.Cells.Item(89, 11).Value = 3682.49994
$ws.Cells.Item(89, 13).Value = 2245.50006
$ws.Cells.Item(131, 8).Value = 15954104
$ws.Cells.Item(131, 9).Value = 12457757
$ws.Cells.Item(131, 10).Value = 17610268
$ws.Cells.Item(131, 11).Value = 37373271
$ws.Cells.Item(131, 12).Value = 52830804
$ws.Cells.Item(131, 13).Value = -37368231
$ws.Cells.Item(131, 14).Value = -52840884
$ws.Cells.Item(137, 8).Value = 59121.45
$ws.Cells.Item(137, 9).Value = 89037.25
$ws.Cells.Item(137, 10).Value = 14247.75
$ws.Cells.Item(137, 11).Value = 267111.75
$ws.Cells.Item(137, 12).Value = 42743.25
$ws.Cells.Item(137, 13).Value = -262011.75
$ws.Cells.Item(137, 14).Value = -52943.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 4337.25
$ws.Cells.Item(43, 9).Value = 4337.25
$ws.Cells.Item(43, 11).Value = 4337.25
$ws.Cells.Item(43, 13).Value = -4186.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1402.9286
$ws.Cells.Item(55, 9).Value = 325.8
$ws.Cells.Item(55, 11).Value = 325.8
$ws.Cells.Item(55, 13).Value = -152.8
$ws.Cells.Item(68, 8).Value = 2081.36
$ws.Cells.Item(68, 9).Value = 2039.7142
$ws.Cells.Item(68, 10).Value = 2300
$ws.Cells.Item(68, 11).Value = 2039.7142
$ws.Cells.Item(68, 12).Value = 2300
$ws.Cells.Item(68, 13).Value = -1290.7142
$ws.Cells.Item(68, 14).Value = -3798
$ws.Cells.Item(71, 8).Value = 2081.36
$ws.Cells.Item(71, 9).Value = 2039.7142
$ws.Cells.Item(71, 10).Value = 2300
$ws.Cells.Item(71, 11).Value = 10198.571
$ws.Cells.Item(71, 12).Value = 11500
$ws.Cells.Item(71, 13).Value = -6454.571
$ws.Cells.Item(71, 14).Value = -18988
$ws.Cells.Item(93, 8).Value = 2132
$ws.Cells.Item(93, 9).Value = 2068
$ws.Cells.Item(93, 10).Value = 2388
$ws.Cells.Item(93, 11).Value = 2068
$ws.Cells.Item(93, 12).Value = 2388
$ws.Cells.Item(93, 13).Value = -820
$ws.Cells.Item(93, 14).Value = -4884

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 47026.91
$ws.Cells.Item(122, 9).Value = 47026.91
$ws.Cells.Item(122, 11).Value = 141080.73
$ws.Cells.Item(122, 13).Value = -138630.73
$ws.Cells.Item(124, 8).Value = 32421.812
$ws.Cells.Item(124, 10).Value = 31249.934
$ws.Cells.Item(124, 12).Value = 31249.934
$ws.Cells.Item(124, 14).Value = -41069.934
$ws.Cells.Item(132, 8).Value = 7887.125
$ws.Cells.Item(132, 9).Value = 8079.6665
$ws.Cells.Item(132, 11).Value = 24238.9995
$ws.Cells.Item(132, 13).Value = -21708.9995
$ws.Cells.Item(136, 8).Value = 7077.352
$ws.Cells.Item(136, 9).Value = 6885.7676
$ws.Cells.Item(136, 10).Value = 7826.273
$ws.Cells.Item(136, 11).Value = 20657.3028
$ws.Cells.Item(136, 12).Value = 23478.819
$ws.Cells.Item(136, 13).Value = -18107.3028
$ws.Cells.Item(136, 14).Value = -28578.819
